# Updated cryptos list on Wed Mar 15 23:39:50 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.520.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.660.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.62%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.37%  "

# Row 6 - USDC
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.22%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3615"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.45%  "

# Row 8 - OKB
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.36"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.88%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3254"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.48%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.91%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06998"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.38%  "

# Row 12 - BinanceUSD
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9993"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.892"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.69%  "

# Row 14 - Solana
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -7.21%  "

# Row 15 - WrappedEther
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.659.45"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.80%  "

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.568"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.63%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001044"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.50%  "

# Row 18 - TRON
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06550"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.46%  "

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.08%  "

# Row 20 - Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -9.33%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.930"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.44%  "

# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.68"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.58%  "

# Row 23 - Cosmos
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.02%  "

# Row 24 - WrappedBTC
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.481.14"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.89%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.469"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.77%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.329"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -15.78%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.59%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.52%  "

# Row 29 - WrappedliquidstakedEther2.0
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.842.46"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.73%  "

# Row 30 - was BitcoinCash, now ImmutableX
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.186"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31 - was ImmutableX, now BitcoinCash
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.60"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.89%  "

# Row 32 - HuobiToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.994"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.33%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.648"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -16.91%  "

# Row 34 - WEMIXTOKEN
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.714"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.95%  "

# Row 35 - Stellar
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08383"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.17%  "

# Row 36 - Aptos
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.32"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -9.99%  "

# Row 37 - InternetComputer(DFINITY)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.184"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.38%  "

# Row 38 - Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06051"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.17%  "

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02203"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.95%  "

# Row 40 - was TrustWalletToken, now FraxShare
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.238"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.36%  "

# Row 41 - was FraxShare, now Algorand
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2056"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.80%  "

# Row 42 - was Algorand, now TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.204"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.69%  "

# Row 43 - Frax
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44 - TheSandbox
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5910"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.34%  "

# Row 45 - PancakeSwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.741"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.90%  "

# Row 46 - EnergySwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.59"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.78%  "

# Row 47 - Decentraland
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5598"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.61%  "

# Row 48 - Quant
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.86%  "

# Row 49 - NEARProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.936"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.71%  "

# Row 50 - Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06893"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.37%  "

# Row 51 - Aave
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.19"
